# Update timestamps in column A (rows 2-97) to reflect the new forecast date (2024-08-29)
$timestamps = @(
    "2024-08-29 00:00:00+00:00",
    "2024-08-29 00:15:00+00:00",
    "2024-08-29 00:30:00+00:00",
    "2024-08-29 00:45:00+00:00",
    "2024-08-29 01:00:00+00:00",
    "2024-08-29 01:15:00+00:00",
    "2024-08-29 01:30:00+00:00",
    "2024-08-29 01:45:00+00:00",
    "2024-08-29 02:00:00+00:00",
    "2024-08-29 02:15:00+00:00",
    "2024-08-29 02:30:00+00:00",
    "2024-08-29 02:45:00+00:00",
    "2024-08-29 03:00:00+00:00",
    "2024-08-29 03:15:00+00:00",
    "2024-08-29 03:30:00+00:00",
    "2024-08-29 03:45:00+00:00",
    "2024-08-29 04:00:00+00:00",
    "2024-08-29 04:15:00+00:00",
    "2024-08-29 04:30:00+00:00",
    "2024-08-29 04:45:00+00:00",
    "2024-08-29 05:00:00+00:00",
    "2024-08-29 05:15:00+00:00",
    "2024-08-29 05:30:00+00:00",
    "2024-08-29 05:45:00+00:00",
    "2024-08-29 06:00:00+00:00",
    "2024-08-29 06:15:00+00:00",
    "2024-08-29 06:30:00+00:00",
    "2024-08-29 06:45:00+00:00",
    "2024-08-29 07:00:00+00:00",
    "2024-08-29 07:15:00+00:00",
    "2024-08-29 07:30:00+00:00",
    "2024-08-29 07:45:00+00:00",
    "2024-08-29 08:00:00+00:00",
    "2024-08-29 08:15:00+00:00",
    "2024-08-29 08:30:00+00:00",
    "2024-08-29 08:45:00+00:00",
    "2024-08-29 09:00:00+00:00",
    "2024-08-29 09:15:00+00:00",
    "2024-08-29 09:30:00+00:00",
    "2024-08-29 09:45:00+00:00",
    "2024-08-29 10:00:00+00:00",
    "2024-08-29 10:15:00+00:00",
    "2024-08-29 10:30:00+00:00",
    "2024-08-29 10:45:00+00:00",
    "2024-08-29 11:00:00+00:00",
    "2024-08-29 11:15:00+00:00",
    "2024-08-29 11:30:00+00:00",
    "2024-08-29 11:45:00+00:00",
    "2024-08-29 12:00:00+00:00",
    "2024-08-29 12:15:00+00:00",
    "2024-08-29 12:30:00+00:00",
    "2024-08-29 12:45:00+00:00",
    "2024-08-29 13:00:00+00:00",
    "2024-08-29 13:15:00+00:00",
    "2024-08-29 13:30:00+00:00",
    "2024-08-29 13:45:00+00:00",
    "2024-08-29 14:00:00+00:00",
    "2024-08-29 14:15:00+00:00",
    "2024-08-29 14:30:00+00:00",
    "2024-08-29 14:45:00+00:00",
    "2024-08-29 15:00:00+00:00",
    "2024-08-29 15:15:00+00:00",
    "2024-08-29 15:30:00+00:00",
    "2024-08-29 15:45:00+00:00",
    "2024-08-29 16:00:00+00:00",
    "2024-08-29 16:15:00+00:00",
    "2024-08-29 16:30:00+00:00",
    "2024-08-29 16:45:00+00:00",
    "2024-08-29 17:00:00+00:00",
    "2024-08-29 17:15:00+00:00",
    "2024-08-29 17:30:00+00:00",
    "2024-08-29 17:45:00+00:00",
    "2024-08-29 18:00:00+00:00",
    "2024-08-29 18:15:00+00:00",
    "2024-08-29 18:30:00+00:00",
    "2024-08-29 18:45:00+00:00",
    "2024-08-29 19:00:00+00:00",
    "2024-08-29 19:15:00+00:00",
    "2024-08-29 19:30:00+00:00",
    "2024-08-29 19:45:00+00:00",
    "2024-08-29 20:00:00+00:00",
    "2024-08-29 20:15:00+00:00",
    "2024-08-29 20:30:00+00:00",
    "2024-08-29 20:45:00+00:00",
    "2024-08-29 21:00:00+00:00",
    "2024-08-29 21:15:00+00:00",
    "2024-08-29 21:30:00+00:00",
    "2024-08-29 21:45:00+00:00",
    "2024-08-29 22:00:00+00:00",
    "2024-08-29 22:15:00+00:00",
    "2024-08-29 22:30:00+00:00",
    "2024-08-29 22:45:00+00:00",
    "2024-08-29 23:00:00+00:00",
    "2024-08-29 23:15:00+00:00",
    "2024-08-29 23:30:00+00:00",
    "2024-08-29 23:45:00+00:00"
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $timestamps[$i]
}

# Update the forecasted power / energy columns (B:F) for rows 26-46 with the
# newly recomputed values for the 2024-08-29 dataset
$rowData = @{
    26 = @(0, 0, 0, 0, 0)
    27 = @(0, 0, 0, 0, 0)
    28 = @(0, 0, 0, 0, 0)
    29 = @(1765.557335535685, 0.0017655573355356, 0, 0.0008827786677678, 0.0002206946669419)
    30 = @(28658.85710652669, 0.0286588571065266, 0.0017655573355356, 0.0152122072210311, 0.0038030518052577)
    31 = @(56037.27966308594, 0.0560372796630859, 0.0286588571065266, 0.0423480683848063, 0.0105870170962015)
    32 = @(48097.56840006511, 0.0480975684000651, 0.0560372796630859, 0.0520674240315754, 0.0130168560078938)
    33 = @(83071.54130045572, 0.0830715413004556, 0.0480975684000651, 0.0655845548502603, 0.016396138712565)
    34 = @(119456.5695800781, 0.1194565695800781, 0.0830715413004556, 0.1012640554402668, 0.0253160138600666)
    35 = @(194109.4351399739, 0.1941094351399739, 0.1194565695800781, 0.1567830023600259, 0.0391957505900064)
    36 = @(603478.037923177, 0.6034780379231771, 0.1941094351399739, 0.3987937365315755, 0.09969843413289381)
    37 = @(406300.1061197916, 0.4063001061197916, 0.6034780379231771, 0.5048890720214844, 0.126222268005371)
    38 = @(486426.8177083333, 0.4864268177083333, 0.4063001061197916, 0.4463634619140624, 0.1115908654785156)
    39 = @(588672.568359375, 0.588672568359375, 0.4864268177083333, 0.5375496930338541, 0.1343874232584635)
    40 = @(658117.0885416666, 0.6581170885416666, 0.588672568359375, 0.6233948284505209, 0.1558487071126302)
    41 = @(657008.0100911459, 0.6570080100911458, 0.6581170885416666, 0.6575625493164062, 0.1643906373291015)
    42 = @(532883.3001302084, 0.5328833001302083, 0.6570080100911458, 0.5949456551106771, 0.1487364137776692)
    43 = @(0, 0, 0.5328833001302083, 0.2664416500651042, 0.066610412516276)
    44 = @(0, 0, 0, 0, 0)
    45 = @(0, 0, 0, 0, 0)
    46 = @(0, 0, 0, 0, 0)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item([int]$row, 2).Value = $vals[0]
    $ws.Cells.Item([int]$row, 3).Value = $vals[1]
    $ws.Cells.Item([int]$row, 4).Value = $vals[2]
    $ws.Cells.Item([int]$row, 5).Value = $vals[3]
    $ws.Cells.Item([int]$row, 6).Value = $vals[4]
}

Write-Host "Updated timestamps and forecast values"
